# Crowdin re-import: fill in newly-translated Korean (column I) and
# Czech (column P) strings that were previously blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Korean (column I)
$ws.Range("I2").Value = "정상적으로 로딩되었습니다.${nl}좋은 하루 보내세요!"
$ws.Range("I3").Value = "상황실"
$ws.Range("I4").Value = "최신 버전 업데이트가 있습니다."
$ws.Range("I5").Value = "버전 업데이트가 가능합니다!"
$ws.Range("I6").Value = "현재 버전${nl}최신 버전${nl}"
$ws.Range("I8").Value = "상황종료- 코드 4"
$ws.Range("I10").Value = "가까이 오세요!"
$ws.Range("I12").Value = "남성"
$ws.Range("I13").Value = "여성"

# Czech (column P)
$ws.Range("P8").Value = "Jste ~g~pod kodem 4~s~.${nl}Neni potreba zadna dalsi jednotka."
$ws.Range("P10").Value = "Prilis daleko.${nl}Prosim, priblizte se."
$ws.Range("P12").Value = "Muz"
$ws.Range("P13").Value = "Zena"
